$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'pun'
$ws.Range('B2').Value = 'play'
$ws.Range('A3').Value = 'confusedsupuerhwomani'
$ws.Range('B3').Value = 'cnn'
$ws.Range('D3').Value = 'noodle'
$ws.Range('A4').Value = 'pianist'
$ws.Range('B4').Value = 'team'
$ws.Range('D4').Value = 'pork'
$ws.Range('A5').Value = 'rosehipskarmamile'
$ws.Range('B5').Value = 'player'
$ws.Range('D5').Value = 'rice'
$ws.Range('E5').Value = 'capsule'
$ws.Range('A6').Value = 'beatz'
$ws.Range('B6').Value = 'ad'
$ws.Range('D6').Value = 'one'
$ws.Range('E6').Value = 'osirisrex'
$ws.Range('A7').Value = 'darjeeling'
$ws.Range('B7').Value = 'nt'
$ws.Range('D7').Value = 'best'
$ws.Range('E7').Value = 'spacecraft'
$ws.Range('A8').Value = 'ringtone'
$ws.Range('B8').Value = 'year'
$ws.Range('E8').Value = 'mission'
$ws.Range('A9').Value = 'songwriter'
$ws.Range('B9').Value = 'game'
$ws.Range('D9').Value = 'chinese'
$ws.Range('E9').Value = 'space'
$ws.Range('B10').Value = 'us'
$ws.Range('E10').Value = 'bennu'
$ws.Range('A11').Value = 'ak'
$ws.Range('B11').Value = 'video'
$ws.Range('E11').Value = 'system'
$ws.Range('A12').Value = 'examiner'
$ws.Range('B12').Value = 'league'
$ws.Range('C12').Value = 'country'
$ws.Range('D12').Value = 'make'
$ws.Range('E12').Value = 'solar'
$ws.Range('B13').Value = 'season'
$ws.Range('D13').Value = 'soup'
$ws.Range('E13').Value = 'nasa'
$ws.Range('A14').Value = 'bee'
$ws.Range('B14').Value = 'get'
$ws.Range('C14').Value = 'company'
$ws.Range('D14').Value = 'serve'
$ws.Range('E14').Value = 'collect'
$ws.Range('A15').Value = 'tmz'
$ws.Range('C15').Value = 'us'
$ws.Range('D15').Value = 'duck'
$ws.Range('E15').Value = 'material'
$ws.Range('B16').Value = 'asteroid'
$ws.Range('C16').Value = 'log'
$ws.Range('D16').Value = 'recipe'
$ws.Range('E16').Value = 'rock'
$ws.Range('A17').Value = 'earl'
$ws.Range('B17').Value = 'sample'
$ws.Range('C17').Value = 'business'
$ws.Range('E17').Value = 'scientist'
$ws.Range('A18').Value = 'honey'
$ws.Range('B18').Value = 'feedback'
$ws.Range('C18').Value = 'chip'
$ws.Range('D18').Value = 'dan'
$ws.Range('E18').Value = 'planet'
$ws.Range('A19').Value = 'affairlovin'
$ws.Range('C19').Value = 'world'
$ws.Range('D19').Value = 'steam'
$ws.Range('E19').Value = 'mile'
$ws.Range('A20').Value = 'empuerh'
$ws.Range('B20').Value = 'export'
$ws.Range('C20').Value = 'tv'
$ws.Range('D20').Value = 'sweet'
$ws.Range('E20').Value = 'deliver'
$ws.Range('A21').Value = 'mindgreenage'
$ws.Range('B21').Value = 'last'
$ws.Range('C21').Value = 'gallium'
$ws.Range('E21').Value = 'orbit'
